# Daily data update: insert a new day's row (2024-05-11) at the top of the
# Sheet2 data table, shifting the existing rows (and the running-total /
# notes rows below them) down by one, then refresh the SUM formulas and the
# free-text "notes" cells in column H to reflect the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert a new blank row above the first data row (row 2); this pushes the
# existing data rows, the totals row and the H-column notes down by one.
$ws.Rows.Item(2).Insert()

# Copy the date-format style from the (now shifted) row below so the new
# date cell matches the rest of the column, then fill in the new day's data.
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(2, 1).PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value = 45423
$ws.Cells.Item(2, 2).Value = 203
$ws.Cells.Item(2, 3).Value = 26
$ws.Cells.Item(2, 4).Value = 9
$ws.Cells.Item(2, 5).Value = 168

# The totals row (now row 13) needs its SUM ranges extended to include the
# new row 2.
$ws.Range("B13").Formula = "=SUM(B2:B12)"
$ws.Range("C13:E13").Formula = "=SUM(C2:C12)"

# Update the free-text notes in column H (now shifted to rows 12-16) with
# the refreshed counts.
$ws.Cells.Item(12, 8).Value = "  count: 203,"
$ws.Cells.Item(14, 8).Value = "      'Night hour arrivals': 26,"
$ws.Cells.Item(15, 8).Value = "      'Regular arrivals': 168,"
$ws.Cells.Item(16, 8).Value = "      'Shoulder hour flights': 9"

# Match the author's final selection on this sheet.
$null = $ws.Range("E16").Select()
